$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Helper: copy all the visual formatting (number format, font, fill,
# alignment) from $src cell to $dst cell, reusing the workbook's existing
# style records instead of synthesising new ones.
# ---------------------------------------------------------------------------
function Copy-Style($src, $dst) {
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# ===========================================================================
# Table "Tabela: Tasks" (F2:K8) - fill in tasktype_id / activity_id columns
# (G, H) that existed but were blank, and extend rows 7-8 with the same
# G/H/J/K columns that rows 4-6 already have.
# ===========================================================================

# Rows 4-6: G/H cells already exist (blank, correctly styled) - just set values.
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1

$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 2

$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 2

# Row 7: add G7, H7 (plain numbers like F7) and J7, K7 (created_at / deleted_at)
$ws.Range("G7").Value = 1
Copy-Style $ws.Range("F7") $ws.Range("G7")

$ws.Range("H7").Value = 1
Copy-Style $ws.Range("F7") $ws.Range("H7")

$ws.Range("J7").Value = 44923.333333333336
Copy-Style $ws.Range("J4") $ws.Range("J7")

$ws.Range("K7").Value = "null"
Copy-Style $ws.Range("K4") $ws.Range("K7")

# Row 8: add G8, H8 (plain numbers like F8) and J8, K8 (created_at / deleted_at)
$ws.Range("G8").Value = 2
Copy-Style $ws.Range("F8") $ws.Range("G8")

$ws.Range("H8").Value = 1
Copy-Style $ws.Range("F8") $ws.Range("H8")

$ws.Range("J8").Value = 44924.333333333336
Copy-Style $ws.Range("J4") $ws.Range("J8")

$ws.Range("K8").Value = "null"
Copy-Style $ws.Range("K4") $ws.Range("K8")

# ===========================================================================
# Table "Tabela: Activities" (A11:H15) - insert a new "total_duration"
# column between "ended_at" and "created_at" (i.e. new column F), shifting
# the former F (created_at) and G (deleted_at) contents one column right.
# ===========================================================================

# --- Header row 12: shift deleted_at (G12) -> H12, created_at (F12) -> G12,
#     then turn F12 into the new "total_duration" header. Work right-to-left
#     so we never overwrite data before it has been copied onward.
$ws.Range("H12").Value = $ws.Range("G12").Value()
Copy-Style $ws.Range("G12") $ws.Range("H12")

$ws.Range("G12").Value = $ws.Range("F12").Value()
# style on G12 is already correct (header label style), left untouched

$ws.Range("F12").Value = "total_duration"
# style on F12 is already correct (header label style), left untouched

# --- Data rows 13-15: shift deleted_at (G) -> H, created_at (F) -> G (still
#     a date), then put the new total_duration time value into F.
$rows = @(
    @{ Row = 13; Total = 0.74305555555555547 },
    @{ Row = 14; Total = 0.22916666666666666 },
    @{ Row = 15; Total = 0.35416666666666669 }
)

foreach ($item in $rows) {
    $r = $item.Row

    $g = $ws.Range("G$r")
    $f = $ws.Range("F$r")
    $h = $ws.Range("H$r")

    # deleted_at ("null") moves from G to H, keeping its plain text style.
    $h.Value = $g.Value()
    Copy-Style $g $h

    # created_at (date) moves from F to G, keeping the date style.
    $g.Value = $f.Value()
    Copy-Style $f $g

    # total_duration (new time value) goes into F, styled like time_elapsed.
    $f.Value = $item.Total
    Copy-Style $ws.Range("I4") $f
}

# --- Grow the section merge from A11:G11 to A11:H11, formatting the new
#     H11 cell like the rest of the merged label row.
$ws.Range("A11:G11").UnMerge()
$ws.Range("H11").Value = ""
Copy-Style $ws.Range("A11") $ws.Range("H11")
$ws.Range("A11:H11").Merge()

# ===========================================================================
# Column widths: column G (now holding total_duration / activity_id values)
# should match the bestFit width already used by columns C-F.
# ===========================================================================
$ws.Range("C1:G1").EntireColumn.AutoFit()
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth()

# ===========================================================================
# Final selection, matching the saved workbook's cursor position.
# ===========================================================================
$ws.Range("H20").Select()
